# Auto-generated: update cryptocurrency Price (D) and Volume(1h) (E) columns
# on Sheet1 to refreshed values from the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.82"
$ws.Range("E2").Value = "'-6.37%"
$ws.Range("D3").Value = "'39.94"
$ws.Range("E3").Value = "'-9.98%"
$ws.Range("E4").Value = "'-5.56%"
$ws.Range("D5").Value = "'0.07771"
$ws.Range("E5").Value = "'-7.15%"
$ws.Range("D6").Value = "'4.320"
$ws.Range("E6").Value = "'-2.29%"
$ws.Range("D7").Value = "'1.644"
$ws.Range("E7").Value = "'-15.21%"
$ws.Range("D8").Value = "'0.9199"
$ws.Range("E8").Value = "'-5.44%"
$ws.Range("D9").Value = "'0.09796"
$ws.Range("E9").Value = "'-13.45%"
$ws.Range("D10").Value = "'0.1737"
$ws.Range("E10").Value = "'-8.68%"
$ws.Range("D11").Value = "'0.08931"
$ws.Range("E11").Value = "'-8.04%"
$ws.Range("D12").Value = "'0.04394"
$ws.Range("E12").Value = "'-4.60%"
$ws.Range("D13").Value = "'7.049"
$ws.Range("E13").Value = "'-15.05%"
$ws.Range("D14").Value = "'0.1059"
$ws.Range("E14").Value = "'-0.25%"
$ws.Range("D15").Value = "'0.001249"
$ws.Range("E15").Value = "'-3.66%"
$ws.Range("E16").Value = "'-3.51%"
$ws.Range("E18").Value = "'1.99%"
$ws.Range("D19").Value = "'0.3370"
$ws.Range("E19").Value = "'0.34%"
$ws.Range("D20").Value = "'0.1363"
$ws.Range("E20").Value = "'-1.79%"
$ws.Range("E21").Value = "'0.19%"
$ws.Range("D22").Value = "'0.04151"
$ws.Range("E22").Value = "'-0.71%"
$ws.Range("D23").Value = "'0.001207"
$ws.Range("E23").Value = "'-2.63%"
$ws.Range("D24").Value = "'0.004090"
$ws.Range("E24").Value = "'-7.65%"
$ws.Range("D25").Value = "'0.0001226"
$ws.Range("E25").Value = "'-5.71%"
$ws.Range("D26").Value = "'0.0002993"
$ws.Range("E26").Value = "'0.45%"
$ws.Range("D38").Value = "'0.02386"
$ws.Range("E38").Value = "'-12.16%"
$ws.Range("D39").Value = "'0.05175"
$ws.Range("E39").Value = "'-8.16%"
$ws.Range("D40").Value = "'0.007986"
$ws.Range("E40").Value = "'1.95%"
$ws.Range("D41").Value = "'0.1328"
$ws.Range("E41").Value = "'-6.25%"
$ws.Range("D42").Value = "'0.007453"
$ws.Range("E42").Value = "'1.99%"
$ws.Range("D43").Value = "'0.001979"
$ws.Range("E43").Value = "'-2.98%"
$ws.Range("D44").Value = "'0.008065"
$ws.Range("E44").Value = "'-7.29%"
$ws.Range("D45").Value = "'0.3335"
$ws.Range("E45").Value = "'-4.97%"
$ws.Range("D46").Value = "'0.00006724"
$ws.Range("E46").Value = "'-2.80%"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("E47").Value = "'0.45%"
$ws.Range("D48").Value = "'0.003424"
$ws.Range("E48").Value = "'-1.79%"
$ws.Range("E49").Value = "'16.67%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.45%"
$ws.Range("D51").Value = "'0.0002010"
$ws.Range("E51").Value = "'0.45%"
